$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Terminal La Palmera de La Serena - Espárragos" (row 2 and row 3)
# have had their Fecha / Volumen / Precio mínimo / Precio máximo / Precio promedio ponderado /
# Precio $/Kg values swapped between the two rows.

# New row 2 values (previously held by row 3)
$ws.Range("D2").Value = 44875
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = 1650
$ws.Range("P2").Value = 1650

# New row 3 values (previously held by row 2)
$ws.Range("D3").Value = 44547
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = 1550
$ws.Range("P3").Value = 1550
